# Update "想去人数" (column F) counts to match regenerated gh-pages data (commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(6, 6).Value2 = 345  # F6: 343 -> 345
$ws.Cells.Item(7, 6).Value2 = 1148  # F7: 1144 -> 1148
$ws.Cells.Item(8, 6).Value2 = 440  # F8: 438 -> 440
$ws.Cells.Item(9, 6).Value2 = 7047  # F9: 7044 -> 7047
$ws.Cells.Item(10, 6).Value2 = 82  # F10: 81 -> 82
$ws.Cells.Item(13, 6).Value2 = 7937  # F13: 7928 -> 7937
$ws.Cells.Item(16, 6).Value2 = 5484  # F16: 5481 -> 5484
$ws.Cells.Item(18, 6).Value2 = 2385  # F18: 2383 -> 2385
$ws.Cells.Item(19, 6).Value2 = 1010  # F19: 1008 -> 1010
$ws.Cells.Item(20, 6).Value2 = 4552  # F20: 4551 -> 4552
$ws.Cells.Item(21, 6).Value2 = 290  # F21: 289 -> 290
$ws.Cells.Item(25, 6).Value2 = 351  # F25: 350 -> 351
$ws.Cells.Item(26, 6).Value2 = 250  # F26: 249 -> 250
$ws.Cells.Item(28, 6).Value2 = 2245  # F28: 2242 -> 2245
$ws.Cells.Item(30, 6).Value2 = 258  # F30: 257 -> 258
$ws.Cells.Item(32, 6).Value2 = 115  # F32: 108 -> 115
$ws.Cells.Item(34, 6).Value2 = 3  # F34: 1 -> 3
$ws.Cells.Item(36, 6).Value2 = 1467  # F36: 1465 -> 1467
$ws.Cells.Item(39, 6).Value2 = 2253  # F39: 2250 -> 2253
$ws.Cells.Item(40, 6).Value2 = 2203  # F40: 2202 -> 2203

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value2 = 71  # F3: 70 -> 71
$ws.Cells.Item(4, 6).Value2 = 56  # F4: 54 -> 56
$ws.Cells.Item(5, 6).Value2 = 22  # F5: 21 -> 22
$ws.Cells.Item(7, 6).Value2 = 93  # F7: 92 -> 93

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(7, 6).Value2 = 345  # F7: 343 -> 345
$ws.Cells.Item(8, 6).Value2 = 1148  # F8: 1145 -> 1148
$ws.Cells.Item(9, 6).Value2 = 440  # F9: 438 -> 440
$ws.Cells.Item(10, 6).Value2 = 7047  # F10: 7044 -> 7047
$ws.Cells.Item(11, 6).Value2 = 82  # F11: 81 -> 82
$ws.Cells.Item(14, 6).Value2 = 7937  # F14: 7928 -> 7937
$ws.Cells.Item(17, 6).Value2 = 5484  # F17: 5481 -> 5484
$ws.Cells.Item(19, 6).Value2 = 2385  # F19: 2383 -> 2385
$ws.Cells.Item(20, 6).Value2 = 1010  # F20: 1008 -> 1010
$ws.Cells.Item(21, 6).Value2 = 4552  # F21: 4551 -> 4552
$ws.Cells.Item(22, 6).Value2 = 290  # F22: 289 -> 290
$ws.Cells.Item(25, 6).Value2 = 71  # F25: 70 -> 71
$ws.Cells.Item(27, 6).Value2 = 56  # F27: 54 -> 56
$ws.Cells.Item(28, 6).Value2 = 351  # F28: 350 -> 351
$ws.Cells.Item(29, 6).Value2 = 250  # F29: 249 -> 250
$ws.Cells.Item(31, 6).Value2 = 2245  # F31: 2242 -> 2245
$ws.Cells.Item(33, 6).Value2 = 258  # F33: 257 -> 258
$ws.Cells.Item(35, 6).Value2 = 115  # F35: 108 -> 115
$ws.Cells.Item(37, 6).Value2 = 3  # F37: 1 -> 3
$ws.Cells.Item(39, 6).Value2 = 22  # F39: 21 -> 22
$ws.Cells.Item(40, 6).Value2 = 1467  # F40: 1465 -> 1467
$ws.Cells.Item(43, 6).Value2 = 2253  # F43: 2250 -> 2253
$ws.Cells.Item(45, 6).Value2 = 2203  # F45: 2202 -> 2203
$ws.Cells.Item(49, 6).Value2 = 93  # F49: 92 -> 93
